# --- Fixed update to excel issue ---
# 1) Rename the "Requested quantity" header on the existing sheets.
# 2) Add a new "PO Forecast" worksheet with forecasted PO quantities.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the header style (bold, centered, bordered) from an existing sheet.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Reuse the date style from an existing sheet's date column.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)  # xlPasteFormats

$wsForecast.Range("A2").Value = 45354.99999999999
$wsForecast.Range("B2").Value = 0
$wsForecast.Range("C2").Value = -22.62577769687611
$wsForecast.Range("D2").Value = 18.11426139711349
$wsForecast.Range("A3").Value = 45375.99999999999
$wsForecast.Range("B3").Value = 24
$wsForecast.Range("C3").Value = 2.943911315494826
$wsForecast.Range("D3").Value = 44.04116907919953
$wsForecast.Range("A4").Value = 45382.99999999999
$wsForecast.Range("B4").Value = 33
$wsForecast.Range("C4").Value = 11.41296421001746
$wsForecast.Range("D4").Value = 53.68542414480055
$wsForecast.Range("A5").Value = 45389.99999999999
$wsForecast.Range("B5").Value = 41
$wsForecast.Range("C5").Value = 21.01853350113297
$wsForecast.Range("D5").Value = 62.53930541544448
$wsForecast.Range("A6").Value = 45396.99999999999
$wsForecast.Range("B6").Value = 50
$wsForecast.Range("C6").Value = 30.11860504052673
$wsForecast.Range("D6").Value = 70.34435552191738
$wsForecast.Range("A7").Value = 45403.99999999999
$wsForecast.Range("B7").Value = 59
$wsForecast.Range("C7").Value = 39.98635621220967
$wsForecast.Range("D7").Value = 78.75772912265595
$wsForecast.Range("A8").Value = 45410.99999999999
$wsForecast.Range("B8").Value = 68
$wsForecast.Range("C8").Value = 47.27440454712624
$wsForecast.Range("D8").Value = 88.95356142994318
$wsForecast.Range("A9").Value = 45417.99999999999
$wsForecast.Range("B9").Value = 77
$wsForecast.Range("C9").Value = 55.12158181129309
$wsForecast.Range("D9").Value = 97.48155947601174
$wsForecast.Range("A10").Value = 45424.99999999999
$wsForecast.Range("B10").Value = 85
$wsForecast.Range("C10").Value = 64.48941843699059
$wsForecast.Range("D10").Value = 106.5796407741963
$wsForecast.Range("A11").Value = 45431.99999999999
$wsForecast.Range("B11").Value = 94
$wsForecast.Range("C11").Value = 73.61005974840907
$wsForecast.Range("D11").Value = 113.6626633601726
$wsForecast.Range("A12").Value = 45438.99999999999
$wsForecast.Range("B12").Value = 103
$wsForecast.Range("C12").Value = 82.65096235084616
$wsForecast.Range("D12").Value = 122.4590120443446

# Restore the originally active sheet/selection.
$wsWeekly.Activate()
